$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 4180.125
$ws.Range("I53").Value = 425
$ws.Range("K53").Value = 425
$ws.Range("M53").Value = 212
$ws.Range("H112").Value = 3907235.2
$ws.Range("J112").Value = 1034.3
$ws.Range("L112").Value = 3102.9
$ws.Range("N112").Value = -5318.9
$ws.Range("H129").Value = 734.8570999999999
$ws.Range("H137").Value = 1537.4849
$ws.Range("I137").Value = 1194.3529
$ws.Range("J137").Value = 1902.0625
$ws.Range("K137").Value = 3583.0587
$ws.Range("L137").Value = 5706.1875
$ws.Range("M137").Value = -1033.0587
$ws.Range("N137").Value = -10806.1875
$ws.Range("H138").Value = 2603.75
$ws.Range("J138").Value = 3214.2917
$ws.Range("L138").Value = 9642.875100000001
$ws.Range("N138").Value = -19922.8751

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3429.9092
$ws.Range("I32").Value = 2356.3438
$ws.Range("K32").Value = 2356.3438
$ws.Range("M32").Value = -2069.3438
$ws.Range("H61").Value = 3476.4285
$ws.Range("I61").Value = 2917.2104
$ws.Range("J61").Value = 4657
$ws.Range("K61").Value = 2917.2104
$ws.Range("L61").Value = 4657
$ws.Range("M61").Value = -2705.2104
$ws.Range("N61").Value = -5081
$ws.Range("H136").Value = 3476.4285
$ws.Range("I136").Value = 2917.2104
$ws.Range("J136").Value = 4657
$ws.Range("K136").Value = 8751.6312
$ws.Range("L136").Value = 13971
$ws.Range("M136").Value = -6201.6312
$ws.Range("N136").Value = -19071

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 49800
$ws.Range("J59").Value = 49800
$ws.Range("L59").Value = 49800
$ws.Range("N59").Value = -51494
$ws.Range("H134").Value = 2601.4055
$ws.Range("I134").Value = 2656.2354
$ws.Range("J134").Value = 1980
$ws.Range("K134").Value = 7968.706200000001
$ws.Range("L134").Value = 5940
$ws.Range("M134").Value = -5433.706200000001
$ws.Range("N134").Value = -11010

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 23950.227
$ws.Range("I58").Value = 1255.25
$ws.Range("J58").Value = 250900
$ws.Range("K58").Value = 1255.25
$ws.Range("L58").Value = 250900
$ws.Range("M58").Value = -1052.25
$ws.Range("N58").Value = -251306
$ws.Range("H105").Value = 8929983
$ws.Range("I105").Value = 12500736
$ws.Range("J105").Value = 3099.75
$ws.Range("K105").Value = 12500736
$ws.Range("L105").Value = 3099.75
$ws.Range("M105").Value = -12498989
$ws.Range("N105").Value = -6593.75
$ws.Range("H136").Value = 23950.227
$ws.Range("I136").Value = 1255.25
$ws.Range("J136").Value = 250900
$ws.Range("K136").Value = 3765.75
$ws.Range("L136").Value = 752700
$ws.Range("M136").Value = -1215.75
$ws.Range("N136").Value = -757800

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1252
$ws.Range("J48").Value = 1252
$ws.Range("L48").Value = 3756
$ws.Range("N48").Value = -4256
$ws.Range("H68").Value = 1280.9524
$ws.Range("I68").Value = 625
$ws.Range("K68").Value = 1875
$ws.Range("M68").Value = -1064
$ws.Range("H71").Value = 1280.9524
$ws.Range("I71").Value = 625
$ws.Range("K71").Value = 5625
$ws.Range("M71").Value = -1569
$ws.Range("H98").Value = 597.7143
$ws.Range("J98").Value = 648
$ws.Range("L98").Value = 1944
$ws.Range("N98").Value = -4940
$ws.Range("H107").Value = 4010.7932
$ws.Range("I107").Value = 6722.8125
$ws.Range("J107").Value = 672.9231
$ws.Range("K107").Value = 20168.4375
$ws.Range("L107").Value = 2018.7693
$ws.Range("M107").Value = -18248.4375
$ws.Range("N107").Value = -5858.7693
$ws.Range("H131").Value = 818.66
$ws.Range("I131").Value = 409.75
$ws.Range("J131").Value = 835.69794
$ws.Range("K131").Value = 1229.25
$ws.Range("L131").Value = 2507.09382
$ws.Range("M131").Value = 3810.75
$ws.Range("N131").Value = -12587.09382

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2549.8
$ws.Range("I102").Value = 2641.4583
$ws.Range("J102").Value = 2183.1667
$ws.Range("K102").Value = 2641.4583
$ws.Range("L102").Value = 2183.1667
$ws.Range("M102").Value = -1019.4583
$ws.Range("N102").Value = -5427.1667
$ws.Range("H122").Value = 3013.3044
$ws.Range("I122").Value = 2347.25
$ws.Range("J122").Value = 4535.7144
$ws.Range("K122").Value = 7041.75
$ws.Range("L122").Value = 13607.1432
$ws.Range("M122").Value = -4591.75
$ws.Range("N122").Value = -18507.1432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1160
$ws.Range("J22").Value = 2150
$ws.Range("L22").Value = 2150
$ws.Range("N22").Value = -2740
$ws.Range("H27").Value = 1160
$ws.Range("J27").Value = 2150
$ws.Range("L27").Value = 2150
$ws.Range("N27").Value = -2364
$ws.Range("H46").Value = 943.1667
$ws.Range("I46").Value = 552.0769
$ws.Range("J46").Value = 1960
$ws.Range("K46").Value = 552.0769
$ws.Range("L46").Value = 1960
$ws.Range("M46").Value = -364.0769
$ws.Range("N46").Value = -2336
$ws.Range("H136").Value = 75771.42999999999
$ws.Range("I136").Value = 251750
$ws.Range("J136").Value = 5380
$ws.Range("K136").Value = 755250
$ws.Range("L136").Value = 16140
$ws.Range("M136").Value = -752700
$ws.Range("N136").Value = -21240

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2350
$ws.Range("I81").Value = 1900
$ws.Range("J81").Value = 3250
$ws.Range("K81").Value = 3800
$ws.Range("L81").Value = 6500
$ws.Range("M81").Value = -2739
$ws.Range("N81").Value = -8622
$ws.Range("H84").Value = 2350
$ws.Range("I84").Value = 1900
$ws.Range("J84").Value = 3250
$ws.Range("K84").Value = 19000
$ws.Range("L84").Value = 32500
$ws.Range("M84").Value = -13696
$ws.Range("N84").Value = -43108
$ws.Range("H107").Value = 1370.174
$ws.Range("I107").Value = 738.5454999999999
$ws.Range("J107").Value = 1949.1666
$ws.Range("K107").Value = 2215.6365
$ws.Range("L107").Value = 5847.4998
$ws.Range("M107").Value = -295.6364999999996
$ws.Range("N107").Value = -9687.4998
$ws.Range("H122").Value = 2097.3076
$ws.Range("J122").Value = 2415
$ws.Range("L122").Value = 7245
$ws.Range("N122").Value = -12145
$ws.Range("H132").Value = 2133.739
$ws.Range("I132").Value = 1944.1428
$ws.Range("J132").Value = 2428.6667
$ws.Range("K132").Value = 5832.428400000001
$ws.Range("L132").Value = 7286.000100000001
$ws.Range("M132").Value = -3302.428400000001
$ws.Range("N132").Value = -12346.0001
$ws.Range("H136").Value = 2934433
$ws.Range("I136").Value = 6452333
$ws.Range("K136").Value = 19356999
$ws.Range("M136").Value = -19354449
